# Generate Report for Handback
# The dfc7aa25-...md row has now been handed back (in sync with en-US) for
# both the zh-cn and de-de locales. Update the Overview sheet and the two
# per-locale sheets accordingly, clear the stale "not latest" error detail,
# and record the new handback/handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the dfc7aa25-...md file. Its zh-cn and de-de
# status columns move from "Ready for handoff" to the handed-back state.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the dfc7aa25-...md file.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-09-02 06:56:43"
$zhcn.Range("K3").Value = "2016-09-02 06:56:43"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 12.9

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the dfc7aa25-...md file.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-02 06:56:51"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 12.9
